$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sec Cap")

# Fill in previously-empty PCI-DSS 4.0 (column E) cells with "O" or "X",
# matching the already-populated GDPR 2016/679 (column F) pattern for
# these rows. Existing cell styles are left untouched since we only set
# the Value.
$ws.Range("E42").Value = "O"
$ws.Range("E43").Value = "O"
$ws.Range("E44").Value = "O"
$ws.Range("E46").Value = "O"
$ws.Range("E48").Value = "O"
$ws.Range("E49").Value = "O"
$ws.Range("E50").Value = "O"
$ws.Range("E51").Value = "O"
$ws.Range("E52").Value = "O"
$ws.Range("E53").Value = "O"
$ws.Range("E54").Value = "O"
$ws.Range("E55").Value = "O"
$ws.Range("E56").Value = "O"
$ws.Range("E57").Value = "O"
$ws.Range("E58").Value = "O"
$ws.Range("E59").Value = "O"
$ws.Range("E60").Value = "O"
$ws.Range("E61").Value = "O"
$ws.Range("E64").Value = "O"
$ws.Range("E70").Value = "X"
$ws.Range("E71").Value = "X"
$ws.Range("E77").Value = "X"
$ws.Range("E78").Value = "X"
$ws.Range("E80").Value = "X"
$ws.Range("E82").Value = "O"
$ws.Range("E89").Value = "O"
$ws.Range("E90").Value = "X"
$ws.Range("E91").Value = "X"
$ws.Range("E94").Value = "X"
$ws.Range("E96").Value = "O"
$ws.Range("E97").Value = "O"
$ws.Range("E98").Value = "X"
$ws.Range("E101").Value = "X"
$ws.Range("E102").Value = "X"
$ws.Range("E103").Value = "X"
$ws.Range("E104").Value = "X"
$ws.Range("E105").Value = "X"
$ws.Range("E106").Value = "O"
$ws.Range("E107").Value = "O"
$ws.Range("E108").Value = "X"
$ws.Range("E109").Value = "X"
$ws.Range("E110").Value = "O"
$ws.Range("E111").Value = "X"
$ws.Range("E112").Value = "X"
$ws.Range("E113").Value = "O"
$ws.Range("E116").Value = "X"
$ws.Range("E121").Value = "X"
$ws.Range("E122").Value = "X"
$ws.Range("E123").Value = "X"
$ws.Range("E124").Value = "O"
$ws.Range("E125").Value = "O"
$ws.Range("E126").Value = "X"
$ws.Range("E127").Value = "X"
$ws.Range("E128").Value = "X"
$ws.Range("E129").Value = "X"
$ws.Range("E130").Value = "O"
$ws.Range("E131").Value = "X"
$ws.Range("E132").Value = "X"
$ws.Range("E133").Value = "X"
$ws.Range("E134").Value = "O"
$ws.Range("E135").Value = "O"
$ws.Range("E136").Value = "O"
$ws.Range("E137").Value = "X"
$ws.Range("E138").Value = "X"
$ws.Range("E139").Value = "X"

# Switch the active/selected sheet from "Index" to "Sec Cap" and move the
# selection there, matching the saved view state in the workbook.
$ws.Activate()
$ws.Range("E140").Select()
